# Update "想去人数" (want-to-go count) values in both the 展览 sheet and the
# 全部类型 (all types) aggregate sheet.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1 / name match)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 1215
$wsExpo.Range("F3").Value = 232
$wsExpo.Range("F4").Value = 68
$wsExpo.Range("F7").Value = 5746
$wsExpo.Range("F8").Value = 5084
$wsExpo.Range("F9").Value = 25
$wsExpo.Range("F10").Value = 53
$wsExpo.Range("F12").Value = 59
$wsExpo.Range("F13").Value = 214
$wsExpo.Range("F14").Value = 13

# Sheet "全部类型" (all types aggregate)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 1215
$wsAll.Range("F3").Value = 232
$wsAll.Range("F4").Value = 68
$wsAll.Range("F7").Value = 5746
$wsAll.Range("F8").Value = 5084
$wsAll.Range("F9").Value = 25
$wsAll.Range("F10").Value = 53
$wsAll.Range("F12").Value = 59
$wsAll.Range("F13").Value = 214
$wsAll.Range("F16").Value = 13
